# Update crypto price (D) and volume-change (E) columns to the latest scraped values.
# D-column values are forced to text via a leading apostrophe so Excel does not
# reinterpret numeric-looking strings (e.g. "27.888.63", "0.9967") as numbers,
# matching the original inline-string cell content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'27.888.63"
$ws.Cells.Item(2, 5).Value = "  +1.06%  "
$ws.Cells.Item(3, 4).Value = "'1.772.96"
$ws.Cells.Item(3, 5).Value = "  +0.77%  "
$ws.Cells.Item(4, 4).Value = "'0.9967"
$ws.Cells.Item(4, 5).Value = "  -0.81%  "
$ws.Cells.Item(5, 5).Value = "  -1.16%  "
$ws.Cells.Item(6, 4).Value = "'0.9960"
$ws.Cells.Item(6, 5).Value = "  -0.64%  "
$ws.Cells.Item(7, 4).Value = "'0.4266"
$ws.Cells.Item(7, 5).Value = "  -4.78%  "
$ws.Cells.Item(8, 4).Value = "'0.3610"
$ws.Cells.Item(8, 5).Value = "  -2.74%  "
$ws.Cells.Item(9, 4).Value = "'44.22"
$ws.Cells.Item(9, 5).Value = "  -2.26%  "
$ws.Cells.Item(10, 4).Value = "'0.07466"
$ws.Cells.Item(10, 5).Value = "  -3.00%  "
$ws.Cells.Item(11, 4).Value = "'1.105"
$ws.Cells.Item(11, 5).Value = "  -1.23%  "
$ws.Cells.Item(12, 5).Value = "  -1.11%  "
$ws.Cells.Item(13, 4).Value = "'21.56"
$ws.Cells.Item(13, 5).Value = "  -0.53%  "
$ws.Cells.Item(14, 4).Value = "'6.122"
$ws.Cells.Item(14, 5).Value = "  -0.78%  "
$ws.Cells.Item(15, 5).Value = "  -1.51%  "
$ws.Cells.Item(16, 4).Value = "'1.796.47"
$ws.Cells.Item(16, 5).Value = "  +1.94%  "
$ws.Cells.Item(17, 4).Value = "'91.20"
$ws.Cells.Item(17, 5).Value = "  +0.36%  "
$ws.Cells.Item(18, 4).Value = "'0.00001061"
$ws.Cells.Item(18, 5).Value = "  -1.24%  "
$ws.Cells.Item(19, 5).Value = "  +1.00%  "
$ws.Cells.Item(20, 4).Value = "'0.9964"
$ws.Cells.Item(20, 5).Value = "  -0.68%  "
$ws.Cells.Item(21, 4).Value = "'17.18"
$ws.Cells.Item(21, 5).Value = "  -1.50%  "
$ws.Cells.Item(22, 4).Value = "'5.942"
$ws.Cells.Item(22, 5).Value = "  -3.82%  "
$ws.Cells.Item(23, 4).Value = "'27.886.41"
$ws.Cells.Item(23, 5).Value = "  +0.99%  "
$ws.Cells.Item(24, 5).Value = "  -1.91%  "
$ws.Cells.Item(25, 4).Value = "'2.164"
$ws.Cells.Item(25, 5).Value = "  -6.50%  "
$ws.Cells.Item(26, 4).Value = "'159.96"
$ws.Cells.Item(26, 5).Value = "  +4.29%  "
$ws.Cells.Item(27, 4).Value = "'20.26"
$ws.Cells.Item(27, 5).Value = "  -1.90%  "
$ws.Cells.Item(28, 4).Value = "'1.992.55"
$ws.Cells.Item(28, 5).Value = "  +1.64%  "
$ws.Cells.Item(29, 5).Value = "  -6.05%  "
$ws.Cells.Item(30, 4).Value = "'125.91"
$ws.Cells.Item(30, 5).Value = "  -1.73%  "
$ws.Cells.Item(31, 4).Value = "'1.166"
$ws.Cells.Item(31, 5).Value = "  -1.81%  "
$ws.Cells.Item(32, 4).Value = "'5.691"
$ws.Cells.Item(32, 5).Value = "  -0.81%  "
$ws.Cells.Item(33, 4).Value = "'0.08984"
$ws.Cells.Item(33, 5).Value = "  -2.58%  "
$ws.Cells.Item(34, 5).Value = "  -3.99%  "
$ws.Cells.Item(35, 5).Value = "  +0.04%  "
$ws.Cells.Item(36, 5).Value = "  +0.03%  "
$ws.Cells.Item(37, 4).Value = "'5.068"
$ws.Cells.Item(37, 5).Value = "  +0.26%  "
$ws.Cells.Item(38, 5).Value = "  -2.65%  "
$ws.Cells.Item(39, 4).Value = "'0.6424"
$ws.Cells.Item(39, 5).Value = "  -0.41%  "
$ws.Cells.Item(40, 4).Value = "'0.06051"
$ws.Cells.Item(40, 5).Value = "  -0.94%  "
$ws.Cells.Item(41, 4).Value = "'1.178"
$ws.Cells.Item(41, 5).Value = "  +0.57%  "
$ws.Cells.Item(42, 4).Value = "'0.9955"
$ws.Cells.Item(42, 5).Value = "  -0.72%  "
$ws.Cells.Item(43, 4).Value = "'7.850"
$ws.Cells.Item(43, 5).Value = "  -1.82%  "
$ws.Cells.Item(44, 4).Value = "'1.392"
$ws.Cells.Item(44, 5).Value = "  -0.68%  "
$ws.Cells.Item(45, 5).Value = "  -0.42%  "
$ws.Cells.Item(46, 4).Value = "'0.5963"
$ws.Cells.Item(46, 5).Value = "  -0.32%  "
$ws.Cells.Item(47, 4).Value = "'3.690"
$ws.Cells.Item(47, 5).Value = "  -1.13%  "
$ws.Cells.Item(48, 4).Value = "'124.19"
$ws.Cells.Item(48, 5).Value = "  -1.48%  "
$ws.Cells.Item(49, 5).Value = "  -0.93%  "
$ws.Cells.Item(50, 5).Value = "  +0.64%  "
$ws.Cells.Item(51, 4).Value = "'0.06883"
$ws.Cells.Item(51, 5).Value = "  -0.17%  "
